$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in the header cell D1: "Contol" -> "Control"
$ws.Range("D1").Value = "Control"

# Move the active selection from D7 to D1
$ws.Range("D1").Select()
